# Fix latency units in report sheet:
#  - O2 header: "Utility" -> "Utility (Percent)"
#  - I3:K14 (min/max/average write latency): append " msec" to each value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header O2
$ws.Range("O2").Value = "Utility (Percent)"

# Append " msec" to each latency value in columns I, J, K for rows 3..14
for ($row = 3; $row -le 14; $row++) {
    foreach ($col in @("I", "J", "K")) {
        $cell = $ws.Range("$col$row")
        $current = $cell.Text
        $cell.Value = "$current msec"
    }
}
